$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.149.87"
$ws.Range("E2").Value = "  +2.73%  "

$ws.Range("D3").Value = "2.630.57"
$ws.Range("E3").Value = "  +0.13%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "567.70"
$ws.Range("E5").Value = "  +6.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.35"
$ws.Range("E6").Value = "  +1.91%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.605"
$ws.Range("E8").Value = "  +6.69%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.82"
$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("E10").Value = "  +3.46%  "

$ws.Range("E11").Value = "  +6.20%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.341"
$ws.Range("E12").Value = "  +2.32%  "

$ws.Range("D13").Value = "3.094.50"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").Value = "60.131.98"
$ws.Range("E14").Value = "  +2.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.66"
$ws.Range("E15").Value = "  +3.68%  "

$ws.Range("D16").Value = "2.639.49"
$ws.Range("E16").Value = "  -0.62%  "

$ws.Range("E17").Value = "  +2.73%  "

$ws.Range("E18").Value = "  +4.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "342.37"
$ws.Range("E19").Value = "  +2.39%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.38"
$ws.Range("E20").Value = "  +2.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.30"
$ws.Range("E21").Value = "  +1.08%  "

$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.38"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("E24").Value = "  +5.07%  "

$ws.Range("E25").Value = "  +1.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.30"
$ws.Range("E27").Value = "  +2.37%  "

$ws.Range("E28").Value = "  +4.46%  "

$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  +4.39%  "

$ws.Range("E31").Value = "  +4.78%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "156.98"
$ws.Range("E32").Value = "  +4.15%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.16"
$ws.Range("E33").Value = "  +2.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.08"
$ws.Range("E34").Value = "  +4.79%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.913"
$ws.Range("E35").Value = "  +8.68%  "

$ws.Range("E36").Value = "  +12.36%  "

$ws.Range("E37").Value = "  +5.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "37.48"
$ws.Range("E38").Value = "  +0.88%  "

$ws.Range("E39").Value = "  +6.15%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "302.53"
$ws.Range("E40").Value = "  +7.30%  "

$ws.Range("E41").Value = "  +2.28%  "

$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("E43").Value = "  +0.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0975"
$ws.Range("E44").Value = "  +4.26%  "

$ws.Range("E45").Value = "  +2.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.29"
$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("E47").Value = "  -0.66%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0235"
$ws.Range("E48").Value = "  +5.07%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.76"
$ws.Range("E49").Value = "  +7.04%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "122.87"
$ws.Range("E50").Value = "  +9.71%  "

$ws.Range("D51").Value = "1.955.94"
$ws.Range("E51").Value = "  +0.52%  "
